{"js": "// Lecture7 lab1 exercise4: shorten the \"Figure 6.2 ... \" and\n// \"Figure 6.3 ...\" sentences from \"is a complete and full tree.\" to\n// \"is a complete.\", and move the `_GoBack` \"last edit\" bookmark from the\n// end of the Figure 6.5 paragraph to sit between \"complete\" and the final\n// \".\" in the (now edited) Figure 6.3 paragraph \u2014 exactly mirroring where\n// Word leaves that bookmark after the author's last keystroke.\nconst body = context.document.body;\n\n// Remove the old `_GoBack` bookmark first (currently sitting at the end of\n// the \"Figure 6.5 is full tree.\" paragraph) so only the new placement below\n// survives.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// \"Figure 6.2 is a complete and full tree.\" -> \"Figure 6.2 is a complete.\"\nconst fig62 = body.search(\"Figure 6.2 is a complete and full tree.\", { matchCase: true });\nfig62.load(\"items\");\nawait context.sync();\nfig62.items[0].insertText(\"Figure 6.2 is a complete.\", \"Replace\");\nawait context.sync();\n\n// \"Figure 6.3 is a complete and full tree.\" -> \"Figure 6.3 is a complete.\"\nconst fig63 = body.search(\"Figure 6.3 is a complete and full tree.\", { matchCase: true });\nfig63.load(\"items\");\nawait context.sync();\nfig63.items[0].insertText(\"Figure 6.3 is a complete.\", \"Replace\");\nawait context.sync();\n\n// Re-seat `_GoBack` right after \"Figure 6.3 is a complete\" (i.e. just before\n// the trailing period), splitting that run in two just like the diff shows.\nconst fig63stem = body.search(\"Figure 6.3 is a complete\", { matchCase: true });\nfig63stem.load(\"items\");\nawait context.sync();\nfig63stem.items[0].getRange(\"After\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Lecture7 lab1 exercise4: shorten the \"Figure 6.2 ... \" and\n# \"Figure 6.3 ...\" sentences from \"is a complete and full tree.\" to\n# \"is a complete.\", and move the `_GoBack` \"last edit\" bookmark from the\n# end of the Figure 6.5 paragraph to sit between \"complete\" and the final\n# \".\" in the (now edited) Figure 6.3 paragraph -- exactly mirroring where\n# Word leaves that bookmark after the author's last keystroke.\n\n$d = $word.ActiveDocument\n\n# Wdistinct constants used below (spelled out since we can't rely on the\n# Word PIA's enum being loaded):\n#   wdReplaceAll    = 2\n#   wdCollapseEnd   = 0\n\n# Drop the bookmark's old position (end of the \"Figure 6.5\" paragraph) --\n# it gets re-created at its new spot further down.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# \"Figure 6.2 is a complete and full tree.\" -> \"Figure 6.2 is a complete.\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Figure 6.2 is a complete and full tree.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Figure 6.2 is a complete.\", 2) | Out-Null\n\n# \"Figure 6.3 is a complete and full tree.\" -> \"Figure 6.3 is a complete.\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Figure 6.3 is a complete and full tree.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Figure 6.3 is a complete.\", 2) | Out-Null\n\n# Re-seat `_GoBack` right after \"Figure 6.3 is a complete\" (i.e. just before\n# the trailing period) -- this splits that run in two, same as the diff.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Figure 6.3 is a complete\") | Out-Null\n$rng.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $rng) | Out-Null\n"}
